# Add "PMcoarse" as a new ignored-pollutant name to the
# tbl_pollutant_ignore table on the "pollutant" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pollutant")

# The table's "Insert Row" placeholder is row 2 (A1:A2 already reserved
# for the table). Writing a value there promotes it to a real data row.
$ws.Range("A2").Value = "PMcoarse"

# Leave the newly entered cell selected, matching the state Excel saves
# after a user types a value and the table grows to include it.
$ws.Range("A2").Select() | Out-Null
